$d = $word.ActiveDocument

# --- Hunk 1: paragraph 20 (empty paragraph after "tempo -> time") ---
# old pPr: <w:rPr><w:b/><w:sz w:val="28"/></w:rPr>
# new pPr: pStyle=PargrafodaLista, ind left=1440, rPr sz=24
$p20 = $d.Paragraphs.Item(20)
$p20.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009726AB" w:rsidRDefault="009726AB" w:rsidP="009726AB"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:ind w:left="1440"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr></w:p>')

# --- Hunk 2: paragraph 29 (the paragraph holding the _GoBack bookmark) ---
# Its pPr becomes pStyle=PargrafodaLista / ind left=1440 / rPr b+sz32 (bookmarks removed from it),
# and five new bulleted paragraphs are added after it, with the _GoBack bookmark now trailing
# the final new paragraph ("tempoVW -> timeView").
$p29 = $d.Paragraphs.Item(29)
$p29.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00086110" w:rsidRDefault="00086110" w:rsidP="00086110"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:ind w:left="1440"/><w:rPr><w:b/><w:sz w:val="32"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>Classe: persistence.php</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>Não possui variáveis declaradas</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:sz w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>Classe</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>totalra.php</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:sz w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>crimeVW -&gt; crimeView</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:sz w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>tempoVW -&gt; timeView</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

Write-Output "edit applied"
